$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45180 -> 45181) for every data row (rows 2 through 397).
$ws.Range("C2:C397").Value = 45181
